$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.852.10'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.827.53'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.76%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.87'
$ws.Range('E5').Value = '  -1.01%  '
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4576'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3672'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07165'
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8733'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07809'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.51'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.845.47'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.317'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.351'
$ws.Range('E15').Value = '  -3.06%  '
$ws.Range('E16').Value = '  -4.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008699'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.875.52'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.45'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.977'
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.44'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.993'
$ws.Range('E24').Value = '  +3.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.57'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.12'
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.993'
$ws.Range('E27').Value = '  -3.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.64'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.907'
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08775'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.102'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7390'
$ws.Range('E32').Value = '  -4.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.471'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.127'
$ws.Range('E34').Value = '  -3.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.500'
$ws.Range('E35').Value = '  -6.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.081'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01938'
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05120'
$ws.Range('E38').Value = '  -2.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.909'
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.906'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4942'
$ws.Range('E41').Value = '  -3.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1589'
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.222'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4641'
$ws.Range('E45').Value = '  -3.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.12'
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.11'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.595'
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06062'
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.60'
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.43'
$ws.Range('E51').Value = '  -0.94%  '
